# notes and gui updates
#
# Strike through two completed to-do items in the notes:
#   - "Grid.initializeGrid() - add pits"
#   - "Written- #5 on pdf"
#
# Applying Font.StrikeThrough on the paragraph's Range sets <w:strike/>
# on both the paragraph mark run properties (pPr/rPr) and on every run
# in the paragraph, matching how Word itself records a strikethrough
# that spans an entire paragraph.

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -like "Grid.initializeGrid()*add pits*") {
        $p.Range.Font.StrikeThrough = $true
    }
    elseif ($t -like "Written- #5*on pdf*") {
        $p.Range.Font.StrikeThrough = $true
    }
}
